$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

# --- Add the two new rows of bet data (rows 34 and 35) ---

# Row 34
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = 45226
$ws.Range("C34").Value = 1
$ws.Range("D34").Formula = "=F33"
$ws.Range("E34").Value = 871
$ws.Range("F34").Formula = "=D34+E34"
$ws.Range("G34").Value = "ESPORTS"
$ws.Range("H34").Value = "WORLDS 2023"
$ws.Range("I34").Value = "DK"
$ws.Range("J34").Value = "PRIMER BARON MAPA 1"
$ws.Range("K34").Value = 1
$ws.Range("L34").Value = 0
$ws.Range("M34").Formula = "=ROUND((F34/`$D`$31-1)*100, 3)+`$M`$29"

# Row 35
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = 45226
$ws.Range("C35").Value = 1
$ws.Range("D35").Formula = "=F34"
$ws.Range("E35").Value = 835
$ws.Range("F35").Formula = "=D35+E35"
$ws.Range("G35").Value = "ESPORTS"
$ws.Range("H35").Value = "WORLDS 2023"
$ws.Range("I35").Value = "WBG"
$ws.Range("J35").Value = "GANA 1 MAPA EN LA SERIE"
$ws.Range("K35").Value = 1
$ws.Range("L35").Value = 0
$ws.Range("M35").Formula = "=ROUND((F35/`$D`$31-1)*100, 3)+`$M`$29"

# --- Carry over the same formatting the rest of the table uses ---
$ws.Range("B33").Copy()
$ws.Range("B34:B35").PasteSpecial(-4122)
$ws.Range("J33").Copy()
$ws.Range("J34:J35").PasteSpecial(-4122)
$ws.Range("M33").Copy()
$ws.Range("M34:M35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the selection to match the new active cell on the sheet
$ws.Range("F26").Select() | Out-Null

$wb.Save() | Out-Null
